$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ObjectR")
$ws2 = $wb.Worksheets.Item("TestD")

# New automation-locator rows appended to the "ObjectR" sheet (Author search
# locators). Written in this cell order so the shared-string table grows in
# the same sequence as the source edit.
$ws1.Range("C13").Value = "//input[@title='Search']"
$ws1.Range("A13").Value = "text_search"
$ws1.Range("A13").Font.Bold = $true
$ws1.Range("B13").Value = "XPATH"

$ws1.Range("A14").Value = "btn_search"
$ws1.Range("A14").Font.Bold = $true
$ws1.Range("B14").Value = "XPATH"
$ws1.Range("C14").Value = "//div[@class='FPdoLc VlcLAe']//input[@value='Google Search']"

# Move the active tab / selection: "TestD" was the active sheet before the
# edit, "ObjectR" becomes active afterwards.
[void]$ws2.Range("D17").Select()
[void]$ws1.Range("C17").Select()
